$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 40000
$ws.Range("I7").Value = 40000
$ws.Range("K7").Value = 40000
$ws.Range("M7").Value = -39888
$ws.Range("H9").Value = 357.18182
$ws.Range("I9").Value = 382.5
$ws.Range("J9").Value = 326.8
$ws.Range("K9").Value = 382.5
$ws.Range("L9").Value = 326.8
$ws.Range("M9").Value = -213.5
$ws.Range("N9").Value = -664.8
$ws.Range("H14").Value = 40000
$ws.Range("I14").Value = 40000
$ws.Range("K14").Value = 40000
$ws.Range("M14").Value = -39809
$ws.Range("H28").Value = 421.9
$ws.Range("I28").Value = 357
$ws.Range("J28").Value = 1006
$ws.Range("K28").Value = 357
$ws.Range("L28").Value = 1006
$ws.Range("M28").Value = 128
$ws.Range("N28").Value = -1976
$ws.Range("H33").Value = 1443.6
$ws.Range("I33").Value = 1388.85
$ws.Range("J33").Value = 1662.6
$ws.Range("K33").Value = 1388.85
$ws.Range("L33").Value = 1662.6
$ws.Range("M33").Value = -1159.85
$ws.Range("N33").Value = -2120.6
$ws.Range("H43").Value = 7253.727
$ws.Range("I43").Value = 6221.222
$ws.Range("J43").Value = 11900
$ws.Range("K43").Value = 6221.222
$ws.Range("L43").Value = 11900
$ws.Range("M43").Value = -6152.222
$ws.Range("N43").Value = -12038
$ws.Range("H53").Value = 445.4
$ws.Range("I53").Value = 424.375
$ws.Range("J53").Value = 482.77777
$ws.Range("K53").Value = 424.375
$ws.Range("L53").Value = 482.77777
$ws.Range("M53").Value = 212.625
$ws.Range("N53").Value = -1756.77777
$ws.Range("H62").Value = 10158.8
$ws.Range("I62").Value = 1669
$ws.Range("K62").Value = 1669
$ws.Range("M62").Value = -1045
$ws.Range("H65").Value = 10158.8
$ws.Range("I65").Value = 1669
$ws.Range("K65").Value = 8345
$ws.Range("M65").Value = -5225
$ws.Range("H86").Value = 4622.923
$ws.Range("I86").Value = 2866.3333
$ws.Range("J86").Value = 5149.9
$ws.Range("K86").Value = 2866.3333
$ws.Range("L86").Value = 5149.9
$ws.Range("M86").Value = -1743.3333
$ws.Range("N86").Value = -7395.9
$ws.Range("H87").Value = 129490
$ws.Range("J87").Value = 129490
$ws.Range("L87").Value = 129490
$ws.Range("N87").Value = -131986
$ws.Range("H89").Value = 4622.923
$ws.Range("I89").Value = 2866.3333
$ws.Range("J89").Value = 5149.9
$ws.Range("K89").Value = 14331.6665
$ws.Range("L89").Value = 25749.5
$ws.Range("M89").Value = -8715.666499999999
$ws.Range("N89").Value = -36981.5
$ws.Range("H90").Value = 129490
$ws.Range("J90").Value = 129490
$ws.Range("L90").Value = 388470
$ws.Range("N90").Value = -400950
$ws.Range("H92").Value = 527.125
$ws.Range("I92").Value = 530.2857
$ws.Range("J92").Value = 505
$ws.Range("K92").Value = 530.2857
$ws.Range("L92").Value = 505
$ws.Range("M92").Value = 717.7143
$ws.Range("N92").Value = -3001
$ws.Range("H106").Value = 10512.5625
$ws.Range("I106").Value = 2441
$ws.Range("K106").Value = 2441
$ws.Range("M106").Value = -1810
$ws.Range("H111").Value = 1635
$ws.Range("I111").Value = 912
$ws.Range("J111").Value = 3081
$ws.Range("K111").Value = 2736
$ws.Range("L111").Value = 9243
$ws.Range("M111").Value = 331
$ws.Range("N111").Value = -15377
$ws.Range("H116").Value = 5608.4546
$ws.Range("I116").Value = 4983.846
$ws.Range("J116").Value = 6510.6665
$ws.Range("K116").Value = 4983.846
$ws.Range("L116").Value = 6510.6665
$ws.Range("M116").Value = -1541.846
$ws.Range("N116").Value = -13394.6665
$ws.Range("H118").Value = 380
$ws.Range("I118").Value = 380
$ws.Range("K118").Value = 1140
$ws.Range("M118").Value = 517
$ws.Range("H121").Value = 1111
$ws.Range("J121").Value = 1111
$ws.Range("L121").Value = 3333
$ws.Range("N121").Value = -6827
$ws.Range("H132").Value = 830.9792
$ws.Range("I132").Value = 838.1087
$ws.Range("K132").Value = 2514.3261
$ws.Range("M132").Value = 15.67389999999978
$ws.Range("H135").Value = 312.66666
$ws.Range("I135").Value = 312.66666
$ws.Range("K135").Value = 2813.99994
$ws.Range("M135").Value = -278.9999399999997
$ws.Range("H138").Value = 2387.74
$ws.Range("I138").Value = 1691.45
$ws.Range("J138").Value = 2561.8125
$ws.Range("K138").Value = 5074.35
$ws.Range("L138").Value = 7685.4375
$ws.Range("M138").Value = 65.64999999999964
$ws.Range("N138").Value = -17965.4375
$ws.Range("H141").Value = 8044.5713
$ws.Range("I141").Value = 8735.333000000001
$ws.Range("J141").Value = 3900
$ws.Range("K141").Value = 26205.999
$ws.Range("L141").Value = 11700
$ws.Range("M141").Value = -21025.999
$ws.Range("N141").Value = -22060

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 2759.125
$ws.Range("I4").Value = 2915.2
$ws.Range("J4").Value = 2499
$ws.Range("K4").Value = 2915.2
$ws.Range("L4").Value = 2499
$ws.Range("M4").Value = -2799.2
$ws.Range("N4").Value = -2731
$ws.Range("H26").Value = 2080
$ws.Range("I26").Value = 2080
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 2080
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -1750
$ws.Range("H32").Value = 5460202.5
$ws.Range("I32").Value = 7830494
$ws.Range("J32").Value = 42393.57
$ws.Range("K32").Value = 7830494
$ws.Range("L32").Value = 42393.57
$ws.Range("M32").Value = -7830207
$ws.Range("N32").Value = -42967.57
$ws.Range("H45").Value = 31252940
$ws.Range("I45").Value = 41669004
$ws.Range("K45").Value = 41669004
$ws.Range("M45").Value = -41668627
$ws.Range("H46").Value = 4659.4287
$ws.Range("J46").Value = 4745.1665
$ws.Range("L46").Value = 4745.1665
$ws.Range("N46").Value = -5383.1665
$ws.Range("H51").Value = 42000
$ws.Range("J51").Value = 42000
$ws.Range("L51").Value = 42000
$ws.Range("N51").Value = -43512
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").ClearContents()
$ws.Range("N58").Value = 0
$ws.Range("H61").Value = 19234460
$ws.Range("I61").Value = 16132223
$ws.Range("J61").Value = 31255634
$ws.Range("K61").Value = 16132223
$ws.Range("L61").Value = 31255634
$ws.Range("M61").Value = -16132011
$ws.Range("N61").Value = -31256058
$ws.Range("H74").Value = 6256880.5
$ws.Range("I74").Value = 8066985
$ws.Range("K74").Value = 8066985
$ws.Range("M74").Value = -8066111
$ws.Range("H77").Value = 6256880.5
$ws.Range("I77").Value = 8066985
$ws.Range("K77").Value = 40334925
$ws.Range("M77").Value = -40330557
$ws.Range("H88").Value = 2691.3684
$ws.Range("I88").Value = 2177.4167
$ws.Range("J88").Value = 3572.4285
$ws.Range("K88").Value = 2177.4167
$ws.Range("L88").Value = 3572.4285
$ws.Range("M88").Value = -1771.4167
$ws.Range("N88").Value = -4384.4285
$ws.Range("H91").Value = 2691.3684
$ws.Range("I91").Value = 2177.4167
$ws.Range("J91").Value = 3572.4285
$ws.Range("K91").Value = 2177.4167
$ws.Range("L91").Value = 3572.4285
$ws.Range("M91").Value = -773.4167000000002
$ws.Range("N91").Value = -6380.4285
$ws.Range("J97").Value = 4598
$ws.Range("L97").Value = 4598
$ws.Range("N97").Value = -5590
$ws.Range("H122").Value = 3632.8
$ws.Range("I122").Value = 2913.5715
$ws.Range("J122").Value = 4262.125
$ws.Range("K122").Value = 8740.7145
$ws.Range("L122").Value = 12786.375
$ws.Range("M122").Value = -6290.7145
$ws.Range("N122").Value = -17686.375
$ws.Range("H132").Value = 8551767
$ws.Range("I132").Value = 13891848
$ws.Range("K132").Value = 41675544
$ws.Range("M132").Value = -41673014
$ws.Range("H136").Value = 19234460
$ws.Range("I136").Value = 16132223
$ws.Range("J136").Value = 31255634
$ws.Range("K136").Value = 48396669
$ws.Range("L136").Value = 93766902
$ws.Range("M136").Value = -48394119
$ws.Range("N136").Value = -93772002

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2261.7778
$ws.Range("I22").Value = 1209.6666
$ws.Range("K22").Value = 1209.6666
$ws.Range("M22").Value = -1036.6666
$ws.Range("H86").Value = 3376.7368
$ws.Range("I86").Value = 3059
$ws.Range("K86").Value = 3059
$ws.Range("M86").Value = -1936
$ws.Range("H89").Value = 3376.7368
$ws.Range("I89").Value = 3059
$ws.Range("K89").Value = 15295
$ws.Range("M89").Value = -9679
$ws.Range("H94").Value = 1406.2273
$ws.Range("I94").Value = 1767.909
$ws.Range("J94").Value = 1044.5454
$ws.Range("K94").Value = 1767.909
$ws.Range("L94").Value = 1044.5454
$ws.Range("M94").Value = -1316.909
$ws.Range("N94").Value = -1946.5454
$ws.Range("H107").Value = 1837.3704
$ws.Range("I107").Value = 1454.3334
$ws.Range("K107").Value = 1454.3334
$ws.Range("M107").Value = 465.6666
$ws.Range("H134").Value = 2697781.5
$ws.Range("I134").Value = 1520.8049
$ws.Range("J134").Value = 11910005
$ws.Range("K134").Value = 4562.4147
$ws.Range("L134").Value = 35730015
$ws.Range("M134").Value = -2027.4147
$ws.Range("N134").Value = -35735085

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 3861
$ws.Range("I10").Value = 3962.8
$ws.Range("J10").Value = 3691.3333
$ws.Range("K10").Value = 3962.8
$ws.Range("L10").Value = 3691.3333
$ws.Range("M10").Value = -3823.8
$ws.Range("N10").Value = -3969.3333
$ws.Range("H31").Value = 1383641.5
$ws.Range("J31").Value = 2919684
$ws.Range("L31").Value = 2919684
$ws.Range("N31").Value = -2920274
$ws.Range("H34").Value = 1383641.5
$ws.Range("J34").Value = 2919684
$ws.Range("L34").Value = 2919684
$ws.Range("N34").Value = -2920088
$ws.Range("H56").Value = 25000
$ws.Range("J56").Value = 25000
$ws.Range("L56").Value = 25000
$ws.Range("N56").Value = -26690
$ws.Range("H58").Value = 4077.9167
$ws.Range("I58").Value = 1786.4
$ws.Range("J58").Value = 5714.7144
$ws.Range("K58").Value = 1786.4
$ws.Range("L58").Value = 5714.7144
$ws.Range("M58").Value = -1583.4
$ws.Range("N58").Value = -6120.7144
$ws.Range("H86").Value = 40366.895
$ws.Range("I86").Value = 4426.8667
$ws.Range("J86").Value = 81836.16
$ws.Range("K86").Value = 4426.8667
$ws.Range("L86").Value = 81836.16
$ws.Range("M86").Value = -3303.8667
$ws.Range("N86").Value = -84082.16
$ws.Range("H89").Value = 40366.895
$ws.Range("I89").Value = 4426.8667
$ws.Range("J89").Value = 81836.16
$ws.Range("K89").Value = 22134.3335
$ws.Range("L89").Value = 409180.8
$ws.Range("M89").Value = -16518.3335
$ws.Range("N89").Value = -420412.8
$ws.Range("H107").Value = 2260.75
$ws.Range("I107").Value = 2979.6
$ws.Range("K107").Value = 2979.6
$ws.Range("M107").Value = -1059.6
$ws.Range("H122").Value = 1817.4706
$ws.Range("I122").Value = 1801.75
$ws.Range("J122").Value = 1855.2
$ws.Range("K122").Value = 5405.25
$ws.Range("L122").Value = 5565.6
$ws.Range("M122").Value = -2955.25
$ws.Range("N122").Value = -10465.6
$ws.Range("H132").Value = 5207.1143
$ws.Range("I132").Value = 1793.1072
$ws.Range("J132").Value = 18863.143
$ws.Range("K132").Value = 5379.321599999999
$ws.Range("L132").Value = 56589.429
$ws.Range("M132").Value = -2849.321599999999
$ws.Range("N132").Value = -61649.429
$ws.Range("H134").Value = 4107.533
$ws.Range("I134").Value = 2548
$ws.Range("J134").Value = 7226.6
$ws.Range("K134").Value = 7644
$ws.Range("L134").Value = 21679.8
$ws.Range("M134").Value = -5109
$ws.Range("N134").Value = -26749.8
$ws.Range("H136").Value = 4077.9167
$ws.Range("I136").Value = 1786.4
$ws.Range("J136").Value = 5714.7144
$ws.Range("K136").Value = 5359.200000000001
$ws.Range("L136").Value = 17144.1432
$ws.Range("M136").Value = -2809.200000000001
$ws.Range("N136").Value = -22244.1432

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 324.2
$ws.Range("J107").Value = 304
$ws.Range("L107").Value = 912
$ws.Range("N107").Value = -4752
$ws.Range("H131").Value = 14757
$ws.Range("J131").Value = 27375
$ws.Range("L131").Value = 82125
$ws.Range("N131").Value = -92205
$ws.Range("H132").Value = 1555.5294
$ws.Range("J132").Value = 1287.4166
$ws.Range("L132").Value = 11586.7494
$ws.Range("N132").Value = -16646.7494
$ws.Range("H140").Value = 218558.86
$ws.Range("I140").Value = 218558.86
$ws.Range("K140").Value = 655676.58
$ws.Range("M140").Value = -650496.58

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 292249.25
$ws.Range("I19").Value = 434500
$ws.Range("K19").Value = 434500
$ws.Range("M19").Value = -434212
$ws.Range("H53").Value = 44500
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 44500
$ws.Range("K53").Value = 0
$ws.Range("L53").ClearContents()
$ws.Range("M53").Value = 44500
$ws.Range("N53").Value = -45762
$ws.Range("H80").Value = 2998.8333
$ws.Range("J80").Value = 3177.4167
$ws.Range("L80").Value = 3177.4167
$ws.Range("N80").Value = -5173.4167
$ws.Range("H83").Value = 2998.8333
$ws.Range("J83").Value = 3177.4167
$ws.Range("L83").Value = 15887.0835
$ws.Range("N83").Value = -25871.0835
$ws.Range("H102").Value = 5729.303
$ws.Range("I102").Value = 1558.8422
$ws.Range("K102").Value = 1558.8422
$ws.Range("M102").Value = 63.15779999999995
$ws.Range("H107").Value = 1601.625
$ws.Range("I107").Value = 1801.6666
$ws.Range("K107").Value = 1801.6666
$ws.Range("M107").Value = 118.3334
$ws.Range("H122").Value = 3734.2666
$ws.Range("I122").Value = 3712.25
$ws.Range("J122").Value = 3822.3333
$ws.Range("K122").Value = 11136.75
$ws.Range("L122").Value = 11466.9999
$ws.Range("M122").Value = -8686.75
$ws.Range("N122").Value = -16366.9999
$ws.Range("H126").Value = 3604.3
$ws.Range("I126").Value = 2673.8333
$ws.Range("K126").Value = 8021.499899999999
$ws.Range("M126").Value = -5551.499899999999
$ws.Range("H132").Value = 34486300
$ws.Range("I132").Value = 43482160
$ws.Range("J132").Value = 2155.8333
$ws.Range("K132").Value = 130446480
$ws.Range("L132").Value = 6467.499899999999
$ws.Range("M132").Value = -130443950
$ws.Range("N132").Value = -11527.4999

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 14999
$ws.Range("J3").Value = 17499
$ws.Range("L3").Value = 17499
$ws.Range("N3").Value = -17723
$ws.Range("H14").Value = 1000
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 1000
$ws.Range("K14").Value = 0
$ws.Range("L14").ClearContents()
$ws.Range("M14").Value = 1000
$ws.Range("N14").Value = -1344
$ws.Range("H15").Value = 14999
$ws.Range("J15").Value = 17499
$ws.Range("L15").Value = 17499
$ws.Range("N15").Value = -17839
$ws.Range("H16").Value = 887.6
$ws.Range("I16").Value = 744.1739
$ws.Range("J16").Value = 2537
$ws.Range("K16").Value = 744.1739
$ws.Range("L16").Value = 2537
$ws.Range("M16").Value = -574.1739
$ws.Range("N16").Value = -2877
$ws.Range("H50").Value = 32497.666
$ws.Range("I50").Value = 27499
$ws.Range("K50").Value = 27499
$ws.Range("M50").Value = -26862
$ws.Range("H68").Value = 1834.25
$ws.Range("J68").Value = 1832.5
$ws.Range("L68").Value = 1832.5
$ws.Range("N68").Value = -3330.5
$ws.Range("H71").Value = 1834.25
$ws.Range("J71").Value = 1832.5
$ws.Range("L71").Value = 9162.5
$ws.Range("N71").Value = -16650.5
$ws.Range("H93").Value = 37038716
$ws.Range("I93").Value = 62501620
$ws.Range("J93").Value = 1757.2727
$ws.Range("K93").Value = 62501620
$ws.Range("L93").Value = 1757.2727
$ws.Range("M93").Value = -62500372
$ws.Range("N93").Value = -4253.2727
$ws.Range("H132").Value = 10747.333
$ws.Range("I132").Value = 12296.8
$ws.Range("K132").Value = 36890.39999999999
$ws.Range("M132").Value = -34360.39999999999
$ws.Range("H136").Value = 73618.39999999999
$ws.Range("I136").Value = 16199.75
$ws.Range("K136").Value = 48599.25
$ws.Range("M136").Value = -46049.25

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 27779392
$ws.Range("I107").Value = 31251654
$ws.Range("K107").Value = 93754962
$ws.Range("M107").Value = -93753042
$ws.Range("H116").Value = 191000
$ws.Range("J116").Value = 191000
$ws.Range("L116").Value = 191000
$ws.Range("N116").Value = -200178
$ws.Range("H122").Value = 4287.8887
$ws.Range("I122").Value = 4287.8887
$ws.Range("K122").Value = 12863.6661
$ws.Range("M122").Value = -10413.6661
$ws.Range("H129").Value = 99955
$ws.Range("J129").Value = 99955
$ws.Range("L129").Value = 99955
$ws.Range("N129").Value = -109955
$ws.Range("H132").Value = 5563160.5
$ws.Range("I132").Value = 7226.4
$ws.Range("J132").Value = 16675028
$ws.Range("K132").Value = 21679.2
$ws.Range("L132").Value = 50025084
$ws.Range("M132").Value = -19149.2
$ws.Range("N132").Value = -50030144
$ws.Range("H136").Value = 3821.8948
$ws.Range("I136").Value = 3176.0278
$ws.Range("K136").Value = 9528.0834
$ws.Range("M136").Value = -6978.0834
